# 2021-06 New South Wales Outbreak Paths by Area.xlsx
# "Add files via upload" -- refresh of the per-day gradient hex colour codes
# in column B ("Colour Code") of the "Date Colours" table, plus the four
# rows that previously had no Colour Code value at all (rows 49-52 gain a
# value, row 53 gets the value that used to belong to row 48).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Date Colours")

# New gradient values for B2:B53 (row -> hex colour string), in row order.
$colours = @(
    "#fcfdff", "#f9faff", "#f6f8ff", "#f3f6ff", "#f0f3ff", "#edf1ff", "#eaefff", "#e7ecff", "#e4eaff",
    "#e0e8ff", "#dde5ff", "#dae3ff", "#d7e1ff", "#d4dfff", "#d0dcff", "#cddaff", "#cad8ff", "#c7d6ff",
    "#c3d3ff", "#c0d1ff", "#bdcfff", "#b9cdff", "#b6caff", "#b2c8ff", "#afc6ff", "#abc4ff", "#a8c2ff",
    "#a4bfff", "#a0bdff", "#9dbbff", "#99b9ff", "#95b7ff", "#91b5ff", "#8db2ff", "#89b0ff", "#85aeff",
    "#80acff", "#7caaff", "#77a8ff", "#73a6ff", "#6ea3ff", "#69a1ff", "#639fff", "#5e9dff", "#589bff",
    "#5299ff", "#4b97ff", "#4395ff", "#3b93ff", "#3191ff", "#248fff", "#118dff"
)

$startRow = 2
for ($i = 0; $i -lt $colours.Length; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 2).Value = $colours[$i]
}
